# Regenerate orders with updated distance/sizes.
#
# The experiment's distance conditions and the "large" size code were
# renumbered:
#   D51 -> D55
#   D80 -> D86
#   D64 -> D69
#   S30 -> S31   (S20 and S25 are unchanged)
#
# These tokens show up embedded inside lots of text values across the sheet
# (Condition, Filename_Left, Filename_Right, Distance, Size, ...), so walk
# every cell in the used range and rewrite any text value that contains one
# of the old tokens.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

$firstRow = $used.Row
$firstCol = $used.Column

for ($r = $firstRow; $r -le ($firstRow + $rowCount - 1); $r++) {
    for ($c = $firstCol; $c -le ($firstCol + $colCount - 1); $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -is [string]) {
            $newVal = $val.Replace("D51", "D55").Replace("D80", "D86").Replace("D64", "D69").Replace("S30", "S31")
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
